$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": zero out species counts, remove percentage column values ---
$wsRange = $wb.Worksheets.Item("Range Status")

for ($r = 2; $r -le 7; $r++) {
    $wsRange.Cells.Item($r, 2).Value = 0          # column B -> 0
    $wsRange.Cells.Item($r, 3).ClearContents()    # column C -> cell removed
}

# --- Sheet "High Priority break-up": update values, rename row3 label, drop rows 4-5 ---
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")

$wsBreakup.Cells.Item(2, 1).Value = "Trend New"
$wsBreakup.Cells.Item(2, 2).Value = 14
$wsBreakup.Cells.Item(2, 3).Value = 14.6
$wsBreakup.Cells.Item(2, 4).Value = 14
$wsBreakup.Cells.Item(2, 5).Value = 14.6

$wsBreakup.Cells.Item(3, 1).Value = "IUCN"
$wsBreakup.Cells.Item(3, 2).Value = 82
$wsBreakup.Cells.Item(3, 3).Value = 85.40000000000001
$wsBreakup.Cells.Item(3, 4).Value = 82
$wsBreakup.Cells.Item(3, 5).Value = 85.40000000000001

# Remove the now-obsolete "Range" and "IUCN" rows (old rows 4 and 5)
$wsBreakup.Range("A4:E5").Delete()
